$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to store the given string as literal text,
# matching the original (unstyled) inline-string cells in the sheet -
# plain .Value assignment lets Excel re-interpret digit-led strings
# (prices, percentages) as numbers, so we briefly flip the cell to the
# Text number format, write the value, then clear the format again so
# no stray style index is left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '42.749.73'
Set-TextValue $ws.Range("E2") '  -1.93%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.574.06'
Set-TextValue $ws.Range("E3") '  -0.43%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.04%  '

# Row 5
Set-TextValue $ws.Range("D5") '301.81'
Set-TextValue $ws.Range("E5") '  +0.26%  '

# Row 6
Set-TextValue $ws.Range("D6") '96.72'
Set-TextValue $ws.Range("E6") '  +0.11%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.574'
Set-TextValue $ws.Range("E7") '  -1.07%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.01%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.548'
Set-TextValue $ws.Range("E9") '  -2.27%  '

# Row 10
Set-TextValue $ws.Range("D10") '36.37'
Set-TextValue $ws.Range("E10") '  -1.29%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0809'
Set-TextValue $ws.Range("E11") '  -1.16%  '

# Row 12
Set-TextValue $ws.Range("D12") '7.70'
Set-TextValue $ws.Range("E12") '  -1.69%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +6.34%  '

# Row 14
Set-TextValue $ws.Range("D14") '2.617.06'
Set-TextValue $ws.Range("E14") '  +1.35%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.884'
Set-TextValue $ws.Range("E15") '  -0.75%  '

# Row 16
Set-TextValue $ws.Range("D16") '14.30'
Set-TextValue $ws.Range("E16") '  -0.30%  '

# Row 17
Set-TextValue $ws.Range("D17") '42.844.54'
Set-TextValue $ws.Range("E17") '  -1.53%  '

# Row 18
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range("D18") '0.0₃0999'
Set-TextValue $ws.Range("E18") '  +1.61%  '

# Row 19
$ws.Range("B19").Value = 'InternetComputer(DFINITY)'
$ws.Range("C19").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D19") '12.93'
Set-TextValue $ws.Range("E19") '  +4.59%  '

# Row 20
Set-TextValue $ws.Range("E20") '  -0.55%  '

# Row 21
Set-TextValue $ws.Range("D21") '71.95'
Set-TextValue $ws.Range("E21") '  -1.46%  '

# Row 22
Set-TextValue $ws.Range("D22") '253.64'
Set-TextValue $ws.Range("E22") '  -4.47%  '

# Row 23
Set-TextValue $ws.Range("E23") '  +0.67%  '

# Row 24
Set-TextValue $ws.Range("D24") '2.13'
Set-TextValue $ws.Range("E24") '  -3.32%  '

# Row 25
Set-TextValue $ws.Range("D25") '28.92'
Set-TextValue $ws.Range("E25") '  -0.98%  '

# Row 26
Set-TextValue $ws.Range("E26") '  -0.06%  '

# Row 27
Set-TextValue $ws.Range("D27") '10.28'
Set-TextValue $ws.Range("E27") '  -0.16%  '

# Row 28
Set-TextValue $ws.Range("D28") '37.71'
Set-TextValue $ws.Range("E28") '  -0.86%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.10'
Set-TextValue $ws.Range("E29") '  -2.85%  '

# Row 30
Set-TextValue $ws.Range("D30") '6.02'
Set-TextValue $ws.Range("E30") '  -1.48%  '

# Row 31
Set-TextValue $ws.Range("D31") '154.82'
Set-TextValue $ws.Range("E31") '  +1.70%  '

# Row 32
Set-TextValue $ws.Range("D32") '3.42'
Set-TextValue $ws.Range("E32") '  -5.24%  '

# Row 33
Set-TextValue $ws.Range("D33") '2.17'
Set-TextValue $ws.Range("E33") '  -2.00%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -1.74%  '

# Row 35
Set-TextValue $ws.Range("D35") '0.0802'
Set-TextValue $ws.Range("E35") '  -1.73%  '

# Row 36
Set-TextValue $ws.Range("D36") '18.23'
Set-TextValue $ws.Range("E36") '  +8.33%  '

# Row 37
Set-TextValue $ws.Range("E37") '  -3.15%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -0.87%  '

# Row 39
Set-TextValue $ws.Range("D39") '23.04'
Set-TextValue $ws.Range("E39") '  -4.25%  '

# Row 40
Set-TextValue $ws.Range("D40") '2.13'
Set-TextValue $ws.Range("E40") '  +31.91%  '

# Row 41
Set-TextValue $ws.Range("D41") '3.42'
Set-TextValue $ws.Range("E41") '  -5.69%  '

# Row 42
Set-TextValue $ws.Range("D42") '3.89'
Set-TextValue $ws.Range("E42") '  +0.24%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.0310'
Set-TextValue $ws.Range("E43") '  -1.55%  '

# Row 44
Set-TextValue $ws.Range("D44") '2.081.38'
Set-TextValue $ws.Range("E44") '  +2.16%  '

# Row 45
Set-TextValue $ws.Range("E45") '  +0.20%  '

# Row 46
Set-TextValue $ws.Range("E46") '  +1.17%  '

# Row 47
Set-TextValue $ws.Range("D47") '85.43'
Set-TextValue $ws.Range("E47") '  -3.11%  '

# Row 48
Set-TextValue $ws.Range("D48") '76.05'
Set-TextValue $ws.Range("E48") '  +9.72%  '

# Row 49
Set-TextValue $ws.Range("D49") '106.43'
Set-TextValue $ws.Range("E49") '  +0.62%  '

# Row 50
Set-TextValue $ws.Range("D50") '2.825.49'
Set-TextValue $ws.Range("E50") '  -0.41%  '

# Row 51
Set-TextValue $ws.Range("D51") '1.67'
Set-TextValue $ws.Range("E51") '  +1.41%  '
